$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F4").Value = 7247
$ws1.Range("F6").Value = 422
$ws1.Range("F7").Value = 3733
$ws1.Range("F11").Value = 607

# Sheet "全部类型" (sheet4)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F5").Value = 7247
$ws4.Range("F8").Value = 422
$ws4.Range("F9").Value = 3733
$ws4.Range("F13").Value = 607
